$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.767.28"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "2.623.09"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'596.32"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").Value = "'150.57"
$ws.Range("E6").Value = "  +2.61%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Value = "'0.109"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'5.70"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  +3.31%  "
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").Value = "'27.82"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "3.094.32"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "63.582.25"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "'0.0000152"
$ws.Range("E16").Value = "  +3.62%  "
$ws.Range("D17").Value = "2.628.29"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "'12.31"
$ws.Range("E18").Value = "  +6.83%  "
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").Value = "'348.36"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  +2.71%  "
$ws.Range("D24").Value = "'66.42"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("E25").Value = "  +13.39%  "
$ws.Range("E26").Value = "  +1.84%  "
$ws.Range("D27").Value = "'1.67"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "'559.30"
$ws.Range("E28").Value = "  -5.61%  "
$ws.Range("D29").Value = "'8.23"
$ws.Range("E29").Value = "  +3.39%  "
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'2.04"
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("D33").Value = "0.0₃0845"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").Value = "'1.75"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'5.25"
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").Value = "'168.86"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").Value = "'0.410"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'19.45"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.94"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "'166.69"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").Value = "'39.89"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  +3.72%  "
$ws.Range("D45").Value = "'0.0598"
$ws.Range("E45").Value = "  +4.64%  "
$ws.Range("D46").Value = "'21.59"
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("D47").Value = "'0.630"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").Value = "'0.0251"
$ws.Range("E48").Value = "  +1.49%  "
$ws.Range("D49").Value = "'1.99"
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("D50").Value = "'0.0967"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "'19.42"
$ws.Range("E51").Value = "  +3.51%  "
